# Applies the "Fixed update to excel issue" change:
#  1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on the Weekly Quantity sheet.
#  2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on the Monthly Trend sheet.
#  3. Add a new "PO Forecast" sheet (positioned after the existing sheets) with a
#     forecast table (ds / PO_Forecast / yhat_lower / yhat_upper) of 10 rows.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the "Requested quantity" headers ------------------------
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet at the end of the workbook ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row - copy the existing styled header cells (bold, centered, bordered)
# from the Weekly Quantity sheet so the new header picks up the same style,
# then overwrite the text.
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($wsForecast.Range("C1:D1"))

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows (10 weekly forecast points)
$rows = @(
    @(45200.99999999999, 1, 0.9999999986881065, 1.000000001354865),
    @(45228.99999999999, 1, 0.9999999987354367, 1.00000000119248),
    @(45235.99999999999, 1, 0.9999999985155181, 1.000000001452986),
    @(45242.99999999999, 1, 0.9999999975836148, 1.000000001932103),
    @(45249.99999999999, 1, 0.9999999951518663, 1.000000004288891),
    @(45256.99999999999, 1, 0.9999999919824863, 1.00000000626435),
    @(45263.99999999999, 1, 0.9999999884593878, 1.000000009914756),
    @(45270.99999999999, 1, 0.9999999839565311, 1.000000014020563),
    @(45277.99999999999, 1, 0.9999999798016068, 1.000000019602044),
    @(45284.99999999999, 1, 0.999999974762629, 1.000000024840431)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $values = $rows[$i]
    $wsForecast.Cells.Item($r, 1).Value = $values[0]
    $wsForecast.Cells.Item($r, 2).Value = $values[1]
    $wsForecast.Cells.Item($r, 3).Value = $values[2]
    $wsForecast.Cells.Item($r, 4).Value = $values[3]
}

# Match the "ds" column's date/time number format used elsewhere in the workbook.
$wsForecast.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Keep the originally active sheet selected (adding the new sheet would
# otherwise leave "PO Forecast" as the active tab).
$wsWeekly.Activate()

